$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells keep their original text storage (General/text),
# matching the source data which is stored as inline strings, not numbers.

$ws.Range('D2').Value = '33.976.14'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.783.91'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '221.47'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.47'
$ws.Range('E8').Value = '  -4.50%  '
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('E10').Value = '  +4.72%  '
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').Value = '2.040.20'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = '1.780.24'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.52'
$ws.Range('E14').Value = '  -4.99%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.626'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '33.987.74'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.21'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.98'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.72'
$ws.Range('E19').Value = '  -3.44%  '
$ws.Range('D20').Value = '0.0₃0779'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.69'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('E23').Value = '  -3.80%  '
$ws.Range('E24').Value = '  -2.03%  '
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.99'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.69'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.20'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.49'
$ws.Range('E33').Value = '  -2.54%  '
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('D35').Value = '1.401.59'
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.640'
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.934'
$ws.Range('E39').Value = '  +4.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.46'
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('E43').Value = '  +1.80%  '
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0491'
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').Value = '1.938.77'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.54'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.80'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -1.34%  '
